$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 307, shifting rows 307:320 down to 308:321.
$ws.Rows.Item(307).Insert()

# Populate the newly inserted row 307 with the new record's data.
$ws.Cells.Item(307, 1).Value2 = 10
$ws.Cells.Item(307, 2).Value2 = "Vega Modelo de Temuco"
$ws.Cells.Item(307, 3).Value2 = "La Araucanía"
$ws.Cells.Item(307, 4).Value2 = 44610
$ws.Cells.Item(307, 5).Value2 = 9
$ws.Cells.Item(307, 6).Value2 = 100112037
$ws.Cells.Item(307, 7).Value2 = "Cebollín"
$ws.Cells.Item(307, 8).Value2 = "Sin especificar"
$ws.Cells.Item(307, 9).Value2 = "Primera"
$ws.Cells.Item(307, 10).Value2 = 30
$ws.Cells.Item(307, 11).Value2 = 8000
$ws.Cells.Item(307, 12).Value2 = 8000
$ws.Cells.Item(307, 13).Value2 = 8000
$ws.Cells.Item(307, 14).Value2 = "`$/docena de paquetes"
$ws.Cells.Item(307, 15).Value2 = "Provincia de Cautín"
$ws.Cells.Item(307, 16).Value2 = 667
$ws.Cells.Item(307, 17).Value2 = 12
$ws.Cells.Item(307, 18).Value2 = "Hortaliza"

# Make sure the date cell for the new row keeps the same formatting (style) as
# the other date cells in column D.
$ws.Cells.Item(307, 4).NumberFormat = $ws.Cells.Item(308, 4).NumberFormat
